$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "62.211.97"
$ws.Range("E2").Value2 = "  -2.28%  "
$ws.Range("D3").Value2 = "2.999.65"
$ws.Range("E3").Value2 = "  -2.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "583.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "145.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  -6.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.520"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  -2.99%  "
$ws.Range("D9").Value2 = "2.996.88"
$ws.Range("E9").Value2 = "  -2.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.147"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -6.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "5.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  -2.57%  "
$ws.Range("E12").Value2 = "  -2.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.0000226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  -4.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "34.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -6.44%  "
$ws.Range("E15").Value2 = "  +2.45%  "
$ws.Range("D16").Value2 = "3.486.40"
$ws.Range("E16").Value2 = "  -2.75%  "
$ws.Range("D17").Value2 = "62.155.05"
$ws.Range("E17").Value2 = "  -2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "6.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -3.50%  "
$ws.Range("D19").Value2 = "2.993.78"
$ws.Range("E19").Value2 = "  -2.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "456.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -4.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "13.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.677"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -4.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "7.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -2.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "80.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -1.39%  "
$ws.Range("E25").Value2 = "  -7.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "12.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -5.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "10.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -6.59%  "
$ws.Range("E28").Value2 = "  -0.11%  "
$ws.Range("E30").Value2 = "  -3.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "7.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -7.38%  "
$ws.Range("E32").Value2 = "  -6.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "26.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  -1.32%  "
$ws.Range("E34").Value2 = "  -4.43%  "
$ws.Range("E35").Value2 = "  -4.71%  "
$ws.Range("D36").Value2 = "0.0₃0784"
$ws.Range("E36").Value2 = "  -5.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "5.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  -5.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "2.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -6.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "50.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  -1.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "8.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -3.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "2.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  -12.22%  "
$ws.Range("E42").Value2 = "  -0.09%  "
$ws.Range("B43").Value2 = "TheGraph"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.270"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -7.14%  "
$ws.Range("B44").Value2 = "Bittensor"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "380.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -14.25%  "
$ws.Range("B45").Value2 = "VeChain"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.0350"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -3.46%  "
$ws.Range("B46").Value2 = "Maker"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value2 = "2.745.91"
$ws.Range("E46").Value2 = "  -2.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "38.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -4.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "128.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -3.06%  "
$ws.Range("E50").Value2 = "  -1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "23.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -6.59%  "
